# Updates cryptos list values (Price and Volume(1h) columns) on the active sheet.
# For Price cells whose new value is a plain decimal number (e.g. "206.70"),
# force the cell to Text format before assigning so Excel keeps the exact
# string (including trailing zeros) instead of converting it to a float,
# then reset the style so no extra formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.872.28'
$ws.Range('E2').Value = '  -2.04%  '
$ws.Range('D3').Value = '1.566.70'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.70'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.492'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.59%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.03'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.55%  '
$ws.Range('E9').Value = '  -0.53%  '
$ws.Range('E10').Value = '  -1.01%  '
$ws.Range('D12').Value = '1.789.96'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').Value = '1.558.68'
$ws.Range('E13').Value = '  -1.29%  '
$ws.Range('E14').Value = '  -1.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.514'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.71%  '
$ws.Range('D16').Value = '26.866.71'
$ws.Range('E16').Value = '  -2.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.54'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.44%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.40'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '214.95'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('E20').Value = '  -1.90%  '
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.12'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.58%  '
$ws.Range('E24').Value = '  -0.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '153.96'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.00%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.72'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('E29').Value = '  -0.96%  '
$ws.Range('E30').Value = '  -0.74%  '
$ws.Range('E31').Value = '  -3.43%  '
$ws.Range('E32').Value = '  -0.69%  '
$ws.Range('D33').Value = '1.403.62'
$ws.Range('E33').Value = '  +1.65%  '
$ws.Range('E34').Value = '  -1.47%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.52'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.27%  '
$ws.Range('E36').Value = '  -0.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.936'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.64%  '
$ws.Range('E38').Value = '  -2.48%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.527'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.96%  '
$ws.Range('E40').Value = '  -1.14%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.991'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.15%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.79'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.17%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.32'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.25%  '
$ws.Range('E45').Value = '  +0.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '63.30'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.44%  '
$ws.Range('D47').Value = '1.702.44'
$ws.Range('E47').Value = '  -0.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.24'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.60%  '
$ws.Range('D49').Value = '0.0₇0980'
$ws.Range('E49').Value = '  -2.23%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0953'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.33%  '
$ws.Range('E51').Value = '  -0.83%  '
